$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: "Utility" -> "Utility (Percent)"
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to the Read Latency columns (I: min, J: max, K: average)
# for every data row (rows 3 through 23).
for ($row = 3; $row -le 23; $row++) {
    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $orig = $cell.Text
        $cell.Value = "$orig msec"
    }
}
